$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the numeric-looking Price column updates so that
# values (including trailing zeros, e.g. "0.8100") are preserved exactly as
# strings rather than being coerced to floating point numbers.

$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D26","D28","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.86"
$ws.Range("D3").Value = "23.99"
$ws.Range("D4").Value = "5.371"
$ws.Range("D5").Value = "0.05827"
$ws.Range("D6").Value = "6.462"
$ws.Range("D7").Value = "3.347"
$ws.Range("D8").Value = "0.8100"
$ws.Range("D9").Value = "0.9218"
$ws.Range("D10").Value = "0.1411"
$ws.Range("D11").Value = "0.07352"
$ws.Range("D12").Value = "0.03100"
$ws.Range("D13").Value = "0.03034"
$ws.Range("D14").Value = "0.09368"
$ws.Range("D15").Value = "3.860"
$ws.Range("D16").Value = "0.001547"
$ws.Range("D17").Value = "0.04708"
$ws.Range("D18").Value = "0.0005988"
$ws.Range("D19").Value = "0.006156"
$ws.Range("D20").Value = "0.001244"
$ws.Range("D21").Value = "0.004691"
$ws.Range("D26").Value = "0.1320"
$ws.Range("D28").Value = "0.0002349"
$ws.Range("D40").Value = "0.03846"
$ws.Range("D41").Value = "0.006379"
$ws.Range("D42").Value = "0.1066"
$ws.Range("D43").Value = "0.002941"
$ws.Range("D44").Value = "0.008529"
$ws.Range("D45").Value = "0.00005259"
$ws.Range("D47").Value = "0.6858"
$ws.Range("D48").Value = "0.001821"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").Value = "0.0002000"

# Volume(1h) text column updates
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
